$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Insert a new (blank) row at 19 for the 4th EC line item. This
#        shifts the trailing signature-block rows (23/24) down to (24/25). ---
$ws.Rows.Item(19).Insert() | Out-Null

# --- 2. The new row 19 should get the "closing" (bottom-thick-border) look
#        that row 18 currently has, while row 18 becomes a plain "middle"
#        row like rows 16/17. Move the formatting down by copying formats
#        only (values are set explicitly afterwards). ---
$ws.Range("B18:J18").Copy() | Out-Null
$ws.Range("B19:J19").PasteSpecial(-4122) | Out-Null    # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("B17:J17").Copy() | Out-Null
$ws.Range("B18:J18").PasteSpecial(-4122) | Out-Null    # xlPasteFormats
$excel.CutCopyMode = $false

# --- 3. Restore row 18's own data (only the formatting was overwritten
#        above; the period must stay 2508). ---
$ws.Range("B18").Value = "CC"
$ws.Range("C18").Value = "1143346806"
$ws.Range("D18").Value = "MARIA TERESA MARMOL BARBOZA"
$ws.Range("E18").Value = "2508"
$ws.Range("F18").Value = 56940
$ws.Range("G18").Value = 1423500

# --- 4. Fill in the new 4th EC line item (period 2509) on row 19. ---
$ws.Range("B19").Value = "CC"
$ws.Range("C19").Value = "1143346806"
$ws.Range("D19").Value = "MARIA TERESA MARMOL BARBOZA"
$ws.Range("E19").Value = "2509"
$ws.Range("F19").Value = 56940
$ws.Range("G19").Value = 1423500

# --- 5. Update totals now that a 4th period was added ---
$ws.Range("E11").Value = 227760
$ws.Range("F13").Value = 4

# --- 6. Swap the "Novedad de Ingreso"/"Novedad de Retiro" header order ---
$ws.Range("H15").Value = "Novedad de Ingreso"
$ws.Range("I15").Value = "Novedad de Retiro"

# --- 7. Fix the signature block: blank line now sits directly above the
#        "NOMBRE.../FIRMA..." labels (shifted down one row by the insert). ---
$ws.Range("B24").Value = "___________________________________"
$ws.Range("H24").Value = "___________________________________"

$ws.Range("B25").Value = "NOMBRE DEL REPRESENTANTE LEGAL"
$ws.Range("H25").Value = "FIRMA DEL REPRESENTANTE LEGAL"
